# Apply crypto price/volume refresh (GitHub Actions data pull).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells
# (t="inlineStr" in the source workbook), so every write below is a string
# assignment. Values that look numeric (e.g. "501.90") get a leading
# apostrophe so Excel stores them as text instead of silently converting
# them to a float and dropping the trailing zero / precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.037.50'
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("D3").Value = '2.364.22'
$ws.Range("E3").Value = '  -4.07%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''501.90'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").Value = '''129.54'
$ws.Range("E6").Value = '  -3.43%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -2.82%  '
$ws.Range("D9").Value = '2.367.52'
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = '''0.0985'
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("E12").Value = '  +3.15%  '
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '2.784.14'
$ws.Range("E14").Value = '  -3.96%  '
$ws.Range("D15").Value = '56.029.11'
$ws.Range("E15").Value = '  -3.42%  '
$ws.Range("D16").Value = '''21.43'
$ws.Range("E16").Value = '  -2.69%  '
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("D18").Value = '2.342.43'
$ws.Range("E18").Value = '  -5.31%  '
$ws.Range("D19").Value = '''10.01'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("D21").Value = '''307.51'
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '''65.85'
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = '''0.998'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = '''0.371'
$ws.Range("E26").Value = '  -2.87%  '
$ws.Range("E27").Value = '  -6.23%  '
$ws.Range("D28").Value = '''7.22'
$ws.Range("E28").Value = '  -5.34%  '
$ws.Range("D29").Value = '''172.42'
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("E30").Value = '  -3.93%  '
$ws.Range("E31").Value = '  -3.10%  '
$ws.Range("D33").Value = '''5.78'
$ws.Range("E33").Value = '  -6.28%  '
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").Value = '  -5.90%  '
$ws.Range("D36").Value = '''17.62'
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("E37").Value = '  -6.31%  '
$ws.Range("E38").Value = '  -5.08%  '
$ws.Range("D39").Value = '''36.20'
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("D40").Value = '''0.802'
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("E41").Value = '  -6.22%  '
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("D43").Value = '''129.05'
$ws.Range("E43").Value = '  -5.85%  '
$ws.Range("D44").Value = '''4.69'
$ws.Range("E44").Value = '  -4.92%  '
$ws.Range("E45").Value = '  -2.75%  '
$ws.Range("D46").Value = '''0.0902'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("D47").Value = '''238.50'
$ws.Range("E47").Value = '  -6.93%  '
$ws.Range("E48").Value = '  -2.64%  '
$ws.Range("E49").Value = '  -4.04%  '
$ws.Range("E50").Value = '  -2.57%  '
